# Add a "Public" column (AH) that flags whether a row's dates may be
# shown on the public frontend (1 = public, 0 = concealed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold header style) from the last existing
# header cell (AG1) onto the new header cell (AH1), then set its text.
$ws.Range("AG1").Copy()
$ws.Cells.Item(1, 34).PasteSpecial(-4122)
$ws.Cells.Item(1, 34).Value = "Public"

# Row 2 (Poz-102430) -> public
$ws.Cells.Item(2, 34).Value = 1
# Row 3 (Poz-102431) -> public
$ws.Cells.Item(3, 34).Value = 1
# Row 4 (UGAMS-10541) -> public
$ws.Cells.Item(4, 34).Value = 1
# Row 5 (Bln-1244) -> not public (concealed on the frontend)
$ws.Cells.Item(5, 34).Value = 0

$excel.CutCopyMode = 0
